$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.096.16"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "1.790.22"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.23"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "2.048.36"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.51"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.82%  "
$ws.Range("D14").Value = "1.783.25"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.083.15"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.622"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.18"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.80"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "0.0₃0780"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.29"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0520"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.52%  "
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.62"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.29%  "
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("D35").Value = "1.430.98"
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.643"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.35"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.91%  "
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.56"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.20%  "
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.921"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +7.85%  "
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("E46").Value = "  +3.87%  "
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.64"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "1.949.21"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  +0.07%  "
